# Update symbol list prices/volumes per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "295.83"; "E2"  = "1.94%"
    "D3"  = "31.12";  "E3"  = "1.11%"
    "D4"  = "4.961";  "E4"  = "0.25%"
    "D5"  = "0.07437"; "E5" = "4.28%"
    "D6"  = "2.146";  "E6"  = "18.43%"
    "D7"  = "7.761";  "E7"  = "1.10%"
    "D8"  = "3.750";  "E8"  = "0.43%"
    "D9"  = "0.9151"; "E9"  = "1.84%"
    "D10" = "0.08796"; "E10" = "16.32%"
    "D11" = "0.1707"; "E11" = "3.81%"
    "D12" = "0.08333"; "E12" = "3.78%"
    "D13" = "0.03152"; "E13" = "3.39%"
    "E14" = "0.76%"
    "D15" = "0.001517"; "E15" = "1.43%"
    "D16" = "0.005807"; "E16" = "2.50%"
    "D17" = "3.510";  "E17" = "0.89%"
    "D18" = "2.076";  "E18" = "-1.16%"
    "E19" = "1.66%"
    "D20" = "0.1294"; "E20" = "1.21%"
    "D21" = "3.979";  "E21" = "-1.61%"
    "D22" = "0.2131"; "E22" = "3.93%"
    "D23" = "0.04561"; "E23" = "0.95%"
    "D24" = "0.001210"; "E24" = "-0.19%"
    "D25" = "0.004627"; "E25" = "15.92%"
    "D26" = "0.0001300"; "E26" = "4.10%"
    "D27" = "0.0003397"
    "D39" = "0.01612"; "E39" = "-0.18%"
    "D40" = "0.04481"; "E40" = "2.86%"
    "D41" = "0.007288"; "E41" = "-0.82%"
    "D42" = "0.008982"
    "D43" = "0.1330"; "E43" = "2.06%"
    "D44" = "0.001963"; "E44" = "-1.76%"
    "D45" = "0.009393"; "E45" = "3.91%"
    "D46" = "0.00006117"; "E46" = "2.00%"
    "E47" = "0.10%"
    "E48" = "-0.20%"
    "E49" = "-3.32%"
    "D50" = "0.00002101"; "E50" = "0.10%"
    "D51" = "0.0002001"; "E51" = "0.10%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
